# Updates the cryptos list with latest price / volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.852.78'
$ws.Range("E2").Value = '  +2.34%  '
$ws.Range("D3").Value = '3.820.65'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '662.77'
$ws.Range("E5").Value = '  +6.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.88'
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("D7").Value = '3.822.09'
$ws.Range("E7").Value = '  +1.36%  '
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  +1.67%  '
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.461'
$ws.Range("E11").Value = '  +2.20%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.97'
$ws.Range("E12").Value = '  +5.30%  '
$ws.Range("E13").Value = '  -1.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.71'
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").Value = '4.455.40'
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").Value = '3.802.99'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '70.755.62'
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.83'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.17'
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  +1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.43'
$ws.Range("E21").Value = '  +8.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '478.43'
$ws.Range("E22").Value = '  +2.32%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.714'
$ws.Range("E23").Value = '  +1.77%  '
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.05'
$ws.Range("E25").Value = '  +0.05%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.23'
$ws.Range("E26").Value = '  +1.60%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.37'
$ws.Range("E27").Value = '  +3.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.13'
$ws.Range("E28").Value = '  -1.20%  '
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("D30").Value = '3.970.88'
$ws.Range("E30").Value = '  +1.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.84'
$ws.Range("E31").Value = '  +6.68%  '
$ws.Range("E32").Value = '  +4.04%  '
$ws.Range("E33").Value = '  +2.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.42'
$ws.Range("E34").Value = '  +2.03%  '
$ws.Range("E35").Value = '  +17.29%  '
$ws.Range("B36").Value = 'RenzoRestakedETH'
$ws.Range("C36").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D36").Value = '3.776.16'
$ws.Range("E36").Value = '  +1.32%  '
$ws.Range("B37").Value = 'Binance-PegBSC-USD'
$ws.Range("C37").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  -0.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '9.07'
$ws.Range("E38").Value = '  +1.03%  '
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.42'
$ws.Range("E40").Value = '  +2.69%  '
$ws.Range("E41").Value = '  +2.61%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.971'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.11'
$ws.Range("E44").Value = '  +9.81%  '
$ws.Range("B45").Value = 'USDe'
$ws.Range("C45").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.00'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '45.61'
$ws.Range("E46").Value = '  +5.90%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '158.78'
$ws.Range("E47").Value = '  +3.60%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.01'
$ws.Range("E48").Value = '  +2.83%  '
$ws.Range("B49").Value = 'TheGraph'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.302'
$ws.Range("E49").Value = '  +0.72%  '
$ws.Range("B50").Value = 'ONDO'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.43'
$ws.Range("E50").Value = '  +4.26%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.49'
$ws.Range("E51").Value = '  +0.97%  '
